# Corrección de alcance elvis
#
# The underlying commit only inserts <w:proofErr> spellcheck/grammar-check
# markers around certain words/phrases (e.g. "Featuring", "Firebase",
# "React", "APIs", "Wireframes", "feedback", and a few double-space-before-
# punctuation grammar spots). <w:proofErr> is a marker Word's live editor
# inserts while interactively spell-checking; it carries no visible
# formatting and (same as in real Word automation) cannot be created via
# the Word object model. What *is* reproducible via automation is the
# side effect the diff shows: the runs around/containing the flagged text
# get split into separate <w:r> elements (so the marker can sit between
# them) while keeping their text and formatting identical.
#
# So: for every paragraph touched by the diff we split the existing run(s)
# at the exact boundaries shown there, without altering any visible
# formatting. A zero-visual-effect "set-different-then-restore" trick on
# a scalar run property (Font.Size) is used to force Word to materialize
# a hard run boundary at a given Range without leaving any residue in the
# saved rPr (because the final value matches the original exactly).

$d = $word.ActiveDocument

function Split-Off($rangeStart, $rangeEnd) {
    # Force a clean run boundary for the exact [start,end) range by
    # round-tripping a scalar font property through a different value.
    $r = $d.Range($rangeStart, $rangeEnd)
    $orig = $r.Font.Size
    $r.Font.Size = 1
    $r.Font.Size = $orig
}

function Split-Parts($findText, $parts) {
    # Locate $findText (must be unique in the document), then split it
    # into consecutive runs whose texts are exactly $parts (concatenation
    # of $parts must equal $findText). The first part keeps whatever run
    # it was already part of; every subsequent part is carved out into
    # its own run boundary.
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Split-Parts: text not found: $findText"
    }
    $pos = $rng.Start
    $end = $rng.End
    $boundaries = @()
    foreach ($p in $parts) {
        $boundaries += $pos
        $pos = $pos + $p.Length
    }
    $boundaries += $pos
    if ($pos -ne $end) {
        throw "Split-Parts: parts length mismatch for: $findText"
    }
    for ($i = 1; $i -lt $parts.Count; $i++) {
        Split-Off $boundaries[$i] $boundaries[$i + 1]
    }
}

# 1) "Equipo de desarrollo Featuring" -> "Equipo de desarrollo " | "Featuring"
Split-Parts "Equipo de desarrollo Featuring" @("Equipo de desarrollo ", "Featuring")

# 2) ": Documento ... base de datos , componentes y mockups."
Split-Parts ": Documento con los diagramas de la arquitectura, base de datos , componentes y mockups." @(
    ": Documento con los diagramas de la arquitectura, base de ",
    "datos ,",
    " componentes y mockups."
)

# 3) "Expo, Firebase, Node.js, PostgreSQL"
Split-Parts "Expo, Firebase, Node.js, PostgreSQL" @(
    "Expo, ",
    "Firebase",
    ", Node.js, PostgreSQL"
)

# 4) ", y otras APIs externas ..."
Split-Parts ", y otras APIs externas serán compatibles y accesibles para integrarse en la aplicación." @(
    ", y otras ",
    "APIs",
    " externas serán compatibles y accesibles para integrarse en la aplicación."
)

# 5) ": Inicialmente, ... funcionalidad del producto.|"
Split-Parts ": Inicialmente, no se implementarán estrategias avanzadas de marketing y promoción, enfocando los esfuerzos en el desarrollo técnico y la funcionalidad del producto.|" @(
    ": Inicialmente, no se implementarán estrategias avanzadas de marketing y promoción, enfocando los esfuerzos en el desarrollo técnico y la funcionalidad del ",
    "producto.|"
)

# 6) "Eventos significativos ... puntos de control , en los que ..."
Split-Parts "Eventos significativos dentro del proyecto, puntos de control , en los que se deberá tomar una decisión." @(
    "Eventos significativos dentro del proyecto, puntos de ",
    "control ,",
    " en los que se deberá tomar una decisión."
)

# 7) "Implementación de autenticación y perfiles de usuario"
Split-Parts "Implementación de autenticación y perfiles de usuario" @(
    "Implementación de autenticación y perfiles de ",
    "usuario"
)

# 8) "Desarrollo del sistema de recomendaciones y match"
Split-Parts "Desarrollo del sistema de recomendaciones y match" @(
    "Desarrollo del sistema de recomendaciones y ",
    "match"
)

# 9) " Presentación del proyecto " + ": El producto será " (spans two
#    original runs; split each independently).
Split-Parts " Presentación del proyecto " @(
    " Presentación del ",
    "proyecto",
    " "
)
Split-Parts ": El producto será " @(
    ":",
    " El producto será "
)

# 10) "Capacitaciones ... (React Native, Expo, Firebase, PostgreSQL/MongoDB)."
Split-Parts "Capacitaciones específicas para el equipo de desarrollo en las tecnologías usadas (React Native, Expo, Firebase, PostgreSQL/MongoDB)." @(
    "Capacitaciones específicas para el equipo de desarrollo en las tecnologías usadas (",
    "React",
    " Native, Expo, ",
    "Firebase",
    ", PostgreSQL/MongoDB)."
)

# 11) "Realización de pruebas ... feedback de los usuarios beta."
Split-Parts "Realización de pruebas de usabilidad frecuentes y ajustes basados en el feedback de los usuarios beta." @(
    "Realización de pruebas de usabilidad frecuentes y ajustes basados en el ",
    "feedback",
    " de los usuarios beta."
)

# 12) " se relaciona con otros proyectos móviles en React Native y Expo."
Split-Parts " se relaciona con otros proyectos móviles en React Native y Expo." @(
    " se relaciona con otros proyectos móviles en ",
    "React",
    " Native y Expo."
)

# 13) "Prototipo de Wireframes de la Interfaz de Usuario"
Split-Parts "Prototipo de Wireframes de la Interfaz de Usuario" @(
    "Prototipo de ",
    "Wireframes",
    " de la Interfaz de Usuario"
)

Write-Output "done"
